$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 203; this shifts rows 203:265 down to 204:266
# and extends the used range / dimension to row 266 automatically.
$ws.Rows(203).Insert()

# Populate the newly inserted (blank) row 203 with the new data record.
$ws.Cells.Item(203, 1).Value = 10
$ws.Cells.Item(203, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(203, 3).Value = "La Araucanía"
$ws.Cells.Item(203, 4).Value = 44988
$ws.Cells.Item(203, 5).Value = 9
$ws.Cells.Item(203, 6).Value = 100112005
$ws.Cells.Item(203, 7).Value = "Puerro"
$ws.Cells.Item(203, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 30
$ws.Cells.Item(203, 11).Value = 14000
$ws.Cells.Item(203, 12).Value = 14000
$ws.Cells.Item(203, 13).Value = 14000
$ws.Cells.Item(203, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(203, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(203, 16).Value = 1167
$ws.Cells.Item(203, 17).Value = 12
$ws.Cells.Item(203, 18).Value = "Hortaliza"
